# cập nhật File kế hoạch
# - Xóa nội dung ("x") ở cột D (D17:D25) và cột E (E26:E39) của bảng kế hoạch.
# - Cập nhật vị trí đang chọn trên sheet (selection) sang E41.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "x" markers that used to live in column D for rows 17-25
$ws.Range("D17:D25").ClearContents() | Out-Null

# Clear the "x" markers that used to live in column E for rows 26-39
$ws.Range("E26:E39").ClearContents() | Out-Null

# Move the active selection to match the saved view state (E41)
$ws.Range("E41").Select() | Out-Null
